# Update Goblin_Profits leve price/profit figures (scheduled runner sync)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2852.0952
$ws.Range("I40").Value = 1785.4286
$ws.Range("J40").Value = 3385.4285
$ws.Range("K40").Value = 1785.4286
$ws.Range("L40").Value = 3385.4285
$ws.Range("M40").Value = -1610.4286
$ws.Range("N40").Value = -3735.4285
$ws.Range("H111").Value = 1100
$ws.Range("I111").Value = 1100
$ws.Range("K111").Value = 3300
$ws.Range("M111").Value = -233
$ws.Range("H137").Value = 9208.23
$ws.Range("I137").Value = 13788
$ws.Range("J137").Value = 1880.6
$ws.Range("K137").Value = 41364
$ws.Range("L137").Value = 5641.799999999999
$ws.Range("M137").Value = -38814
$ws.Range("N137").Value = -10741.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2199.5
$ws.Range("I2").Value = 1766
$ws.Range("K2").Value = 1766
$ws.Range("M2").Value = -1653
$ws.Range("H32").Value = 5058.9585
$ws.Range("I32").Value = 5018.0435
$ws.Range("J32").Value = 6000
$ws.Range("K32").Value = 5018.0435
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -4731.0435
$ws.Range("N32").Value = -6574
$ws.Range("H61").Value = 2729.5293
$ws.Range("I61").Value = 2729.5293
$ws.Range("K61").Value = 2729.5293
$ws.Range("M61").Value = -2517.5293
$ws.Range("H74").Value = 1857.76
$ws.Range("I74").Value = 1857.76
$ws.Range("K74").Value = 1857.76
$ws.Range("M74").Value = -983.76
$ws.Range("H77").Value = 1857.76
$ws.Range("I77").Value = 1857.76
$ws.Range("K77").Value = 9288.799999999999
$ws.Range("M77").Value = -4920.799999999999
$ws.Range("H116").Value = 2199.5
$ws.Range("I116").Value = 1766
$ws.Range("K116").Value = 1766
$ws.Range("M116").Value = 528
$ws.Range("H122").Value = 2779657.8
$ws.Range("I122").Value = 3368645
$ws.Range("J122").Value = 3004
$ws.Range("K122").Value = 10105935
$ws.Range("L122").Value = 9012
$ws.Range("M122").Value = -10103485
$ws.Range("N122").Value = -13912
$ws.Range("H132").Value = 2784.7334
$ws.Range("I132").Value = 2212.8948
$ws.Range("K132").Value = 6638.6844
$ws.Range("M132").Value = -4108.6844
$ws.Range("H136").Value = 2729.5293
$ws.Range("I136").Value = 2729.5293
$ws.Range("K136").Value = 8188.5879
$ws.Range("M136").Value = -5638.5879

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2199.5
$ws.Range("I3").Value = 1766
$ws.Range("K3").Value = 1766
$ws.Range("M3").Value = -1652
$ws.Range("H107").Value = 3861.0789
$ws.Range("I107").Value = 1575.0869
$ws.Range("J107").Value = 7366.2666
$ws.Range("K107").Value = 1575.0869
$ws.Range("L107").Value = 7366.2666
$ws.Range("M107").Value = 344.9131
$ws.Range("N107").Value = -11206.2666
$ws.Range("H134").Value = 2922.6924
$ws.Range("I134").Value = 2922.6924
$ws.Range("K134").Value = 8768.0772
$ws.Range("M134").Value = -6233.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 3616.739
$ws.Range("I23").Value = 3611.8572
$ws.Range("K23").Value = 10835.5716
$ws.Range("M23").Value = -10600.5716
$ws.Range("H107").Value = 3450.125
$ws.Range("J107").Value = 1266.6666
$ws.Range("L107").Value = 3799.9998
$ws.Range("N107").Value = -7639.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3863.7778
$ws.Range("I132").Value = 4096.75
$ws.Range("K132").Value = 12290.25
$ws.Range("M132").Value = -9760.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5080
$ws.Range("I7").Value = 5080
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5080
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4968
$ws.Range("N7").Value = $null
$ws.Range("H55").Value = 1300
$ws.Range("I55").Value = 408
$ws.Range("K55").Value = 408
$ws.Range("M55").Value = -235
$ws.Range("H68").Value = 9659.846
$ws.Range("I68").Value = 2157.6
$ws.Range("J68").Value = 14348.75
$ws.Range("K68").Value = 2157.6
$ws.Range("L68").Value = 14348.75
$ws.Range("M68").Value = -1408.6
$ws.Range("N68").Value = -15846.75
$ws.Range("H71").Value = 9659.846
$ws.Range("I71").Value = 2157.6
$ws.Range("J71").Value = 14348.75
$ws.Range("K71").Value = 10788
$ws.Range("L71").Value = 71743.75
$ws.Range("M71").Value = -7044
$ws.Range("N71").Value = -79231.75
$ws.Range("H122").Value = 4998
$ws.Range("I122").Value = 4998
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14994
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12544
$ws.Range("N122").Value = $null
$ws.Range("H126").Value = 5080
$ws.Range("I126").Value = 5080
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15240
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12770
$ws.Range("N126").Value = $null
$ws.Range("H132").Value = 3694.7837
$ws.Range("I132").Value = 3230.6538
$ws.Range("J132").Value = 4791.8184
$ws.Range("K132").Value = 9691.9614
$ws.Range("L132").Value = 14375.4552
$ws.Range("M132").Value = -7161.9614
$ws.Range("N132").Value = -19435.4552
$ws.Range("H139").Value = 98538.57000000001
$ws.Range("J139").Value = 98538.57000000001
$ws.Range("L139").Value = 98538.57000000001
$ws.Range("N139").Value = -108818.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = $null
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = $null
$ws.Range("H107").Value = 4613.7896
$ws.Range("I107").Value = 2039.3529
$ws.Range("J107").Value = 26496.5
$ws.Range("K107").Value = 6118.0587
$ws.Range("L107").Value = 79489.5
$ws.Range("M107").Value = -4198.0587
$ws.Range("N107").Value = -83329.5
$ws.Range("H132").Value = 7443.522
$ws.Range("I132").Value = 6169.4614
$ws.Range("J132").Value = 9099.799999999999
$ws.Range("K132").Value = 18508.3842
$ws.Range("L132").Value = 27299.4
$ws.Range("M132").Value = -15978.3842
$ws.Range("N132").Value = -32359.4
